{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Source diff summary:\n//   \"...disminuir los paros por faltantes y asegurar la continuidad...\"\n// becomes\n//   \"...disminuir los paros por faltantes de importados y asegurar la continuidad...\"\n//\n// i.e. the text \" de importados\" is inserted right after \"por faltantes\"\n// (and right before \" y asegurar la continuidad...\") in the paragraph that\n// starts with \"de esta manera se pretende disminuir los paros...\".\n//\n// (The rest of the source diff only wraps the already-present words\n// \"learning\"/\"Forecast\" with Word's <w:proofErr> spell-check markers \u2014\n// those carry no visible text change and Word's automation object model\n// does not expose authoring of <w:proofErr> elements, so there is nothing\n// further to apply there.)\n\nconst body = context.document.body;\n\n// Find the unique anchor phrase that contains the insertion point.\nconst results = body.search(\"por faltantes y asegurar\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Anchor text \"por faltantes y asegurar\" not found.');\n}\n\nconst target = results.items[0];\n\n// Narrow the match down to just \"por faltantes\" so we can insert right\n// after \"faltantes\" (before the following \" y asegurar...\").\nconst faltantesRanges = target.search(\"por faltantes\", { matchCase: true });\nfaltantesRanges.load(\"items\");\nawait context.sync();\n\nconst insertAfter = faltantesRanges.items[0];\ninsertAfter.insertText(\" de importados\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Source diff summary:\n#   \"...disminuir los paros por faltantes y asegurar la continuidad...\"\n# becomes\n#   \"...disminuir los paros por faltantes de importados y asegurar la continuidad...\"\n#\n# i.e. the text \" de importados\" is inserted right after \"por faltantes\"\n# (and right before \" y asegurar la continuidad...\") in the paragraph that\n# starts with \"de esta manera se pretende disminuir los paros...\".\n#\n# Note: the phrase \"por faltantes\" also appears later in the document\n# (\"... por faltantes y por tiempos inactivos ...\"), so we anchor the\n# Find/Replace on the longer, unique phrase \"por faltantes y asegurar\" to\n# make sure only the intended occurrence is touched.\n#\n# (The rest of the source diff only wraps the already-present words\n# \"learning\"/\"Forecast\" with Word's <w:proofErr> spell-check markers \u2014\n# those carry no visible text change and Word's automation object model\n# does not expose authoring of <w:proofErr> elements, so there is nothing\n# further to apply there.)\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"por faltantes y asegurar\"\n$find.Replacement.Text = \"por faltantes de importados y asegurar\"\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n"}
